# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) worksheet gains three new trailing columns:
#   H: date             -> "2013-11-22"   (kept as literal text)
#   I: legislator_name  -> "吳育昇"
#   J: legislator_id    -> 1322
# applied to the header row and every data row of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "吳育昇"
$legislatorId = 1322
$reportDate = "2013-11-22"

# Find the last used data row from column A (row 2 .. lastRow)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header row (row 1)
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Format column H as text first so the ISO date string ("2013-11-22") is
# kept verbatim instead of being auto-converted into a date serial number.
$dateRange = $ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8))
$dateRange.NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
